# Update "想去人数" (F column) values on the 展览 (sheet1) and 全部类型 (sheet4) sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 304
$wsExhibit.Range("F5").Value = 162
$wsExhibit.Range("F6").Value = 176
$wsExhibit.Range("F7").Value = 311
$wsExhibit.Range("F8").Value = 220
$wsExhibit.Range("F9").Value = 2117
$wsExhibit.Range("F10").Value = 367
$wsExhibit.Range("F11").Value = 5125
$wsExhibit.Range("F12").Value = 110

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 304
$wsAll.Range("F6").Value = 162
$wsAll.Range("F7").Value = 176
$wsAll.Range("F8").Value = 311
$wsAll.Range("F9").Value = 220
$wsAll.Range("F12").Value = 2117
$wsAll.Range("F13").Value = 367
$wsAll.Range("F14").Value = 5125
$wsAll.Range("F15").Value = 110
